$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 434
$wsExhibition.Range("F4").Value = 3507
$wsExhibition.Range("F6").Value = 32
$wsExhibition.Range("F7").Value = 175

# Sheet "演出" (Performance) - sheet2
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 122

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 434
$wsAll.Range("F3").Value = 122
$wsAll.Range("F8").Value = 3507
$wsAll.Range("F10").Value = 32
$wsAll.Range("F12").Value = 175
